# Update the cryptocurrency listing with latest scraped prices / volume figures.
# (Updated cryptos list on Mon Jul 15 14:43:03 UTC 2024 with GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that Excel would otherwise auto-convert to a number
# (e.g. "1.00", "23.00") while keeping it stored as plain text, and without
# leaving a residual number-format style applied to the cell.
function Set-TextValue {
    param($cellRef, $value)
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '62.947.85'
$ws.Range('E2').Value = '  +5.36%  '
$ws.Range('D3').Value = '3.362.22'
$ws.Range('E3').Value = '  +5.61%  '
$ws.Range('E4').Value = '  +0.02%  '
Set-TextValue 'D5' '570.83'
$ws.Range('E5').Value = '  +7.02%  '
$ws.Range('E6').Value = '  +5.85%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '3.366.02'
$ws.Range('E8').Value = '  +5.53%  '
$ws.Range('E9').Value = '  -0.61%  '
Set-TextValue 'D10' '7.41'
$ws.Range('E10').Value = '  +1.50%  '
$ws.Range('E11').Value = '  +5.30%  '
Set-TextValue 'D12' '0.441'
$ws.Range('E12').Value = '  +2.79%  '
$ws.Range('D13').Value = '3.936.75'
$ws.Range('E13').Value = '  +5.38%  '
$ws.Range('E14').Value = '  +0.08%  '
Set-TextValue 'D15' '26.96'
$ws.Range('E15').Value = '  +4.56%  '
$ws.Range('E16').Value = '  +4.97%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '62.976.02'
$ws.Range('E17').Value = '  +5.32%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.380.41'
$ws.Range('E18').Value = '  +5.35%  '
Set-TextValue 'D19' '6.33'
$ws.Range('E19').Value = '  +1.90%  '
Set-TextValue 'D20' '13.92'
$ws.Range('E20').Value = '  +6.02%  '
Set-TextValue 'D21' '8.40'
$ws.Range('E21').Value = '  +3.00%  '
Set-TextValue 'D22' '385.08'
$ws.Range('E22').Value = '  +5.35%  '
$ws.Range('E23').Value = '  +0.27%  '
Set-TextValue 'D24' '0.534'
$ws.Range('E24').Value = '  +2.96%  '
Set-TextValue 'D25' '70.55'
$ws.Range('E25').Value = '  +1.55%  '
Set-TextValue 'D26' '9.34'
$ws.Range('E26').Value = '  +6.97%  '
$ws.Range('E27').Value = '  +6.44%  '
$ws.Range('D28').Value = '0.0₃0964'
$ws.Range('E28').Value = '  +10.66%  '
Set-TextValue 'D29' '1.00'
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('E30').Value = '  +6.27%  '
Set-TextValue 'D31' '23.00'
$ws.Range('E31').Value = '  +3.39%  '
$ws.Range('E32').Value = '  +5.19%  '
$ws.Range('E33').Value = '  +9.71%  '
$ws.Range('E34').Value = '  +3.35%  '
$ws.Range('E35').Value = '  +2.52%  '
Set-TextValue 'D37' '157.77'
$ws.Range('E37').Value = '  +1.61%  '
$ws.Range('E38').Value = '  +12.15%  '
Set-TextValue 'D39' '26.88'
$ws.Range('E39').Value = '  +3.70%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '2.848.34'
$ws.Range('E40').Value = '  +2.27%  '
Set-TextValue 'D41' '0.0740'
$ws.Range('E41').Value = '  +5.98%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D42' '0.0327'
$ws.Range('E42').Value = '  +10.88%  '
Set-TextValue 'D43' '40.93'
$ws.Range('E43').Value = '  +3.89%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D44' '0.744'
$ws.Range('E44').Value = '  +4.44%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D45' '4.26'
$ws.Range('E45').Value = '  +1.38%  '
Set-TextValue 'D46' '1.04'
$ws.Range('E46').Value = '  +5.39%  '
$ws.Range('D47').Value = '3.404.56'
$ws.Range('E47').Value = '  +5.53%  '
Set-TextValue 'D48' '21.93'
$ws.Range('E48').Value = '  +7.25%  '
Set-TextValue 'D49' '297.22'
$ws.Range('E49').Value = '  +13.25%  '
$ws.Range('E50').Value = '  -1.82%  '
Set-TextValue 'D51' '6.30'
$ws.Range('E51').Value = '  +2.98%  '
